$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# (the runtime's ColumnWidth setter quantizes to pixel boundaries, so the
# inputs below are chosen to land as close as possible to the target stored
# widths of 26.85546875 / 24.7109375 characters)
$ws.Columns.Item(1).ColumnWidth = 26
$ws.Columns.Item(2).ColumnWidth = 23.833333333333332

# --- Set the new text values for column B (order mirrors the original authoring
# order so the shared-string table comes out the same) ---
$ws.Range("B10").Value = "1.356 [.302, .398]"
$ws.Range("B11").Value = "1.326 [.261, .371]"
$ws.Range("B13").Value = "1.702 [.574, .891]"

$ws.Range("A15").Value = "offdiagE6W15HRCombined.txt"
$ws.Range("B15").Value = "1.443 [.331, .523]"
$ws.Range("E15").Value = "E6W15HRCombined.png"

$ws.Range("B2").Value = "1.142 [.117, .167]"
$ws.Range("B3").Value = "1.177 [.150, .206]"
$ws.Range("B4").Value = "1.253 [.196, .306]"
$ws.Range("B5").Value = "1.368 [.187, .528]"
$ws.Range("B6").Value = "1.182 [.160, .205]"
$ws.Range("B7").Value = "1.201 [.174, .230]"
$ws.Range("B8").Value = "1.225 [.134, .271]"
$ws.Range("B12").Value = "1.382 [.272, .445]"
$ws.Range("B9").Value = "1.484 [.390, .588]"
$ws.Range("B14").Value = "1.105 [.985, .044]"

# --- Updated / new numeric values ---
$ws.Range("C2").Value = 0.31
$ws.Range("C3").Value = 0.44
$ws.Range("C4").Value = 0.586
$ws.Range("C5").Value = 0.764
$ws.Range("C6").Value = 0.386
$ws.Range("C7").Value = 0.492
$ws.Range("C8").Value = 0.621
$ws.Range("C9").Value = 0.791
$ws.Range("C10").Value = 0.635
$ws.Range("C11").Value = 0.64
$ws.Range("C12").Value = 0.743
$ws.Range("C13").Value = 0.883

# --- New row 15 remaining fields ---
$ws.Range("C15").Value = 0.8
$ws.Range("D15").Value = 3121

# --- Selection ---
$ws.Range("C14").Select()
